$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.626469135284424
$ws.Range("B1").Value = 2.209495544433594
$ws.Range("C1").Value = 5.278173446655273
$ws.Range("D1").Value = 1.469803333282471
$ws.Range("E1").Value = 0.7792685627937317
